$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5789666666666667
$ws.Range("H2").Value = 1.7369
$ws.Range("I2").Value = 0.01523705650035473
$ws.Range("J2").Value = 0.01523705650035472
$ws.Range("M2").Value = 0.1112926666666667
$ws.Range("N2").Value = 0.333878
$ws.Range("O2").Value = 0.01397697460904174
$ws.Range("P2").Value = 0.01397697460904174
$ws.Range("Q2").Value = 0.06443474424444445
$ws.Range("R2").Value = 0.5799126982
$ws.Range("S2").Value = 0.0002129679518219923
$ws.Range("T2").Value = 0.0002129679518219923
$ws.Range("G3").Value = 0.5789666666666667
$ws.Range("H3").Value = 1.7369
$ws.Range("I3").Value = 0.01523705650035473
$ws.Range("J3").Value = 0.01523705650035472
$ws.Range("O3").Value = 0.4165551449121381
$ws.Range("P3").Value = 0.4165551449121381
$ws.Range("Q3").Value = 1.920345781322222
$ws.Range("R3").Value = 17.2831120319
$ws.Range("S3").Value = 0.006347074278539699
$ws.Range("T3").Value = 0.006347074278539699
$ws.Range("G4").Value = 0.5789666666666667
$ws.Range("H4").Value = 1.7369
$ws.Range("I4").Value = 0.01523705650035473
$ws.Range("J4").Value = 0.01523705650035472
$ws.Range("O4").Value = 0.5694678804788202
$ws.Range("P4").Value = 0.5694678804788201
$ws.Range("Q4").Value = 2.625283243366667
$ws.Range("R4").Value = 23.6275491903
$ws.Range("S4").Value = 0.008677014269993037
$ws.Range("T4").Value = 0.008677014269993033
$ws.Range("I5").Value = 0.6545086962501954
$ws.Range("J5").Value = 0.6545086962501954
$ws.Range("M5").Value = 0.1112926666666667
$ws.Range("N5").Value = 0.333878
$ws.Range("O5").Value = 0.01397697460904174
$ws.Range("P5").Value = 0.01397697460904174
$ws.Range("Q5").Value = 2.767798389909777
$ws.Range("R5").Value = 24.910185509188
$ws.Range("S5").Value = 0.00914805142888599
$ws.Range("T5").Value = 0.00914805142888599
$ws.Range("I6").Value = 0.6545086962501954
$ws.Range("J6").Value = 0.6545086962501954
$ws.Range("O6").Value = 0.4165551449121381
$ws.Range("P6").Value = 0.4165551449121381
$ws.Range("S6").Value = 0.2726389648127547
$ws.Range("T6").Value = 0.2726389648127547
$ws.Range("I7").Value = 0.6545086962501954
$ws.Range("J7").Value = 0.6545086962501954
$ws.Range("O7").Value = 0.5694678804788202
$ws.Range("P7").Value = 0.5694678804788201
$ws.Range("S7").Value = 0.3727216800085547
$ws.Range("T7").Value = 0.3727216800085547
$ws.Range("I8").Value = 0.33025424724945
$ws.Range("J8").Value = 0.3302542472494499
$ws.Range("M8").Value = 0.1112926666666667
$ws.Range("N8").Value = 0.333878
$ws.Range("O8").Value = 0.01397697460904174
$ws.Range("P8").Value = 0.01397697460904174
$ws.Range("Q8").Value = 1.396585223442889
$ws.Range("R8").Value = 12.569267010986
$ws.Range("S8").Value = 0.004615955228333754
$ws.Range("T8").Value = 0.004615955228333753
$ws.Range("I9").Value = 0.33025424724945
$ws.Range("J9").Value = 0.3302542472494499
$ws.Range("O9").Value = 0.4165551449121381
$ws.Range("P9").Value = 0.4165551449121381
$ws.Range("S9").Value = 0.1375691058208437
$ws.Range("T9").Value = 0.1375691058208437
$ws.Range("I10").Value = 0.33025424724945
$ws.Range("J10").Value = 0.3302542472494499
$ws.Range("O10").Value = 0.5694678804788202
$ws.Range("P10").Value = 0.5694678804788201
$ws.Range("S10").Value = 0.1880691862002725
$ws.Range("T10").Value = 0.1880691862002725
